{"js": "// Replace the three-digit x one-digit multiplication problems with the\n// newly generated set. Most cells are a straight 1:1 text swap; the\n// fifth data row additionally drops its first problem (its other four\n// problems shift left) and gains a new fifth problem, so that row is\n// rewritten as a whole via the table's `values` grid.\n\nconst simpleReplacements = [\n  [\"894\u00d78=7152\", \"411\u00d72=822\"],\n  [\"422\u00d72=844\", \"409\u00d72=818\"],\n  [\"671\u00d72=1342\", \"625\u00d72=1250\"],\n  [\"844\u00d72=1688\", \"288\u00d72=576\"],\n  [\"538\u00d75=2690\", \"862\u00d76=5172\"],\n  [\"393\u00d77=2751\", \"711\u00d78=5688\"],\n  [\"626\u00d77=4382\", \"462\u00d77=3234\"],\n  [\"861\u00d73=2583\", \"225\u00d78=1800\"],\n  [\"876\u00d74=3504\", \"995\u00d72=1990\"],\n  [\"849\u00d74=3396\", \"826\u00d78=6608\"],\n  [\"962\u00d79=8658\", \"237\u00d76=1422\"],\n  [\"250\u00d75=1250\", \"332\u00d77=2324\"],\n  [\"448\u00d73=1344\", \"955\u00d74=3820\"],\n  [\"739\u00d73=2217\", \"338\u00d72=676\"],\n  [\"815\u00d76=4890\", \"123\u00d75=615\"],\n  [\"143\u00d73=429\", \"276\u00d74=1104\"],\n  [\"297\u00d76=1782\", \"345\u00d73=1035\"],\n  [\"167\u00d74=668\", \"576\u00d73=1728\"],\n  [\"215\u00d72=430\", \"736\u00d75=3680\"],\n  [\"434\u00d72=868\", \"160\u00d73=480\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of simpleReplacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Row that loses its first problem (545\u00d77=3815) and gains a new last\n// problem (644\u00d74=2576); the remaining three problems also get new\n// values. Rewrite the whole row through the table's values grid so the\n// cell count stays at five while every value ends up correct.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst targetRow = table.values.findIndex(\n  (row) => row[0] === \"545\u00d77=3815\" && row[1] === \"231\u00d79=2079\"\n);\nif (targetRow === -1) {\n  throw new Error(\"Could not locate the row to rewrite\");\n}\n\nconst newValues = table.values;\nnewValues[targetRow] = [\"231\u00d79=2079\", \"266\u00d76=1596\", \"149\u00d74=596\", \"292\u00d76=1752\", \"644\u00d74=2576\"];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the three-digit x one-digit multiplication problems with the\n# newly generated set. Most cells are a straight 1:1 text swap done with\n# Find/Replace; the fifth data row (table row 5) additionally drops its\n# first problem (the remaining four shift left) and gains a new fifth\n# problem, so that row is rewritten cell-by-cell instead.\n\n$d = $word.ActiveDocument\n\n$simpleReplacements = @(\n    @(\"894\u00d78=7152\", \"411\u00d72=822\"),\n    @(\"422\u00d72=844\", \"409\u00d72=818\"),\n    @(\"671\u00d72=1342\", \"625\u00d72=1250\"),\n    @(\"844\u00d72=1688\", \"288\u00d72=576\"),\n    @(\"538\u00d75=2690\", \"862\u00d76=5172\"),\n    @(\"393\u00d77=2751\", \"711\u00d78=5688\"),\n    @(\"626\u00d77=4382\", \"462\u00d77=3234\"),\n    @(\"861\u00d73=2583\", \"225\u00d78=1800\"),\n    @(\"876\u00d74=3504\", \"995\u00d72=1990\"),\n    @(\"849\u00d74=3396\", \"826\u00d78=6608\"),\n    @(\"962\u00d79=8658\", \"237\u00d76=1422\"),\n    @(\"250\u00d75=1250\", \"332\u00d77=2324\"),\n    @(\"448\u00d73=1344\", \"955\u00d74=3820\"),\n    @(\"739\u00d73=2217\", \"338\u00d72=676\"),\n    @(\"815\u00d76=4890\", \"123\u00d75=615\"),\n    @(\"143\u00d73=429\", \"276\u00d74=1104\"),\n    @(\"297\u00d76=1782\", \"345\u00d73=1035\"),\n    @(\"167\u00d74=668\", \"576\u00d73=1728\"),\n    @(\"215\u00d72=430\", \"736\u00d75=3680\"),\n    @(\"434\u00d72=868\", \"160\u00d73=480\")\n)\n\nforeach ($pair in $simpleReplacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# Row that loses its first problem (545\u00d77=3815) and gains a new last\n# problem (644\u00d74=2576); the remaining three problems also get new\n# values. Rewrite the five cells of that table row directly so the cell\n# count stays at five while every value ends up correct.\n$table = $d.Tables.Item(1)\n$newRowValues = @(\"231\u00d79=2079\", \"266\u00d76=1596\", \"149\u00d74=596\", \"292\u00d76=1752\", \"644\u00d74=2576\")\nfor ($col = 1; $col -le 5; $col++) {\n    $table.Cell(5, $col).Range.Text = $newRowValues[$col - 1]\n}\n"}
